$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.044.14"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "3.037.41"
$ws.Range("E3").Value = "  -1.09%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.87"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.39"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.16%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -1.13%  "

# Row 9
$ws.Range("D9").Value = "3.037.47"
$ws.Range("E9").Value = "  -0.83%  "

# Row 10
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.153"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -1.00%  "

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -1.05%  "

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.71%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -2.18%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.31"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -2.04%  "

# Row 15
$ws.Range("E15").Value = "  +1.55%  "

# Row 16
$ws.Range("D16").Value = "3.545.60"
$ws.Range("E16").Value = "  -0.96%  "

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.14"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -0.74%  "

# Row 18
$ws.Range("D18").Value = "63.045.77"
$ws.Range("E18").Value = "  -0.32%  "

# Row 19
$ws.Range("D19").Value = "3.048.07"
$ws.Range("E19").Value = "  -0.67%  "

# Row 20
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.26"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +1.03%  "

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -2.60%  "

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -1.74%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +2.26%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.96"
$ws.Range("D25").Style = $style

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.73"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -2.22%  "

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.66"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +8.27%  "

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +1.44%  "

# Row 30
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.67"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -0.05%  "

# Row 31
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.10%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.20"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +0.54%  "

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.48"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +1.05%  "

# Row 34
$ws.Range("E34").Value = "  -2.87%  "

# Row 35
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.05"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +0.48%  "

# Row 36
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0814"
$ws.Range("E36").Value = "  -3.35%  "

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.27"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -1.65%  "

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.92"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -2.92%  "

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.22"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +0.57%  "

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.24"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.69%  "

# Row 41
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.37"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +0.33%  "

# Row 42
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "434.09"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -2.21%  "

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.287"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +0.69%  "

# Row 44
$ws.Range("E44").Value = "  +2.95%  "

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0362"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -0.15%  "

# Row 46
$ws.Range("D46").Value = "2.824.55"
$ws.Range("E46").Value = "  +0.61%  "

# Row 47
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.29"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -4.49%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.53"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -1.19%  "

# Row 49
$ws.Range("E49").Value = "  +0.02%  "

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.05"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +0.10%  "

# Row 51
$ws.Range("E51").Value = "  -1.78%  "
